$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in row 8 with the new timing entry
$ws.Range("A8").Formula = "=(7.9654+7.6745+7.82873)/3"
$ws.Range("B8").Formula = "=A8/60"
$ws.Range("C8").Value = "replacing L = np.linalg.inv(R) with the correct formula"

# Move the active selection to C9, matching the diff's sheetView selection
$ws.Range("C9").Select()
